$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit2"
$ws.Range("C2").Value = "Gpc1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.143896
$ws.Range("H2").Value = 0.431688
$ws.Range("I2").Value = 0.02807111181859822
$ws.Range("J2").Value = 0.02807111181859822
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8478306666666667
$ws.Range("N2").Value = 2.543492
$ws.Range("O2").Value = 0.01460351867535248
$ws.Range("P2").Value = 0.01460351867535248
$ws.Range("Q2").Value = 0.1219994416106667
$ws.Range("R2").Value = 1.097994974496
$ws.Range("S2").Value = 0.000409937005680807
$ws.Range("T2").Value = 0.0004099370056808069

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit2"
$ws.Range("C3").Value = "Gpc1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.143896
$ws.Range("H3").Value = 0.431688
$ws.Range("I3").Value = 0.02807111181859822
$ws.Range("J3").Value = 0.02807111181859822
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.020353
$ws.Range("N3").Value = 21.061059
$ws.Range("O3").Value = 0.1209225617494376
$ws.Range("P3").Value = 0.1209225617494376
$ws.Range("Q3").Value = 1.010200715288
$ws.Range("R3").Value = 9.091806437592
$ws.Range("S3").Value = 0.00339443075225981
$ws.Range("T3").Value = 0.00339443075225981

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit2"
$ws.Range("C4").Value = "Gpc1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.143896
$ws.Range("H4").Value = 0.431688
$ws.Range("I4").Value = 0.02807111181859822
$ws.Range("J4").Value = 0.02807111181859822
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.462291666666667
$ws.Range("N4").Value = 4.386875
$ws.Range("O4").Value = 0.02518734518879435
$ws.Range("P4").Value = 0.02518734518879435
$ws.Range("Q4").Value = 0.2104179216666666
$ws.Range("R4").Value = 1.893761295
$ws.Range("S4").Value = 0.0007070367832082782
$ws.Range("T4").Value = 0.0007070367832082782

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Slit2"
$ws.Range("C5").Value = "Gpc1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.143896
$ws.Range("H5").Value = 0.431688
$ws.Range("I5").Value = 0.02807111181859822
$ws.Range("J5").Value = 0.02807111181859822
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 48.72612633333333
$ws.Range("N5").Value = 146.178379
$ws.Range("O5").Value = 0.8392865743864156
$ws.Range("P5").Value = 0.8392865743864156
$ws.Range("Q5").Value = 7.011494674861333
$ws.Range("R5").Value = 63.103452073752
$ws.Range("S5").Value = 0.02355970727744933
$ws.Range("T5").Value = 0.02355970727744933

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit2"
$ws.Range("C6").Value = "Gpc1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.277274333333334
$ws.Range("H6").Value = 12.831823
$ws.Range("I6").Value = 0.8344071140950421
$ws.Range("J6").Value = 0.8344071140950421
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8478306666666667
$ws.Range("N6").Value = 2.543492
$ws.Range("O6").Value = 0.01460351867535248
$ws.Range("P6").Value = 0.01460351867535248
$ws.Range("Q6").Value = 3.626404349546223
$ws.Range("R6").Value = 32.637639145916
$ws.Range("S6").Value = 0.01218527987353392
$ws.Range("T6").Value = 0.01218527987353392

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit2"
$ws.Range("C7").Value = "Gpc1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.277274333333334
$ws.Range("H7").Value = 12.831823
$ws.Range("I7").Value = 0.8344071140950421
$ws.Range("J7").Value = 0.8344071140950421
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.020353
$ws.Range("N7").Value = 21.061059
$ws.Range("O7").Value = 0.1209225617494376
$ws.Range("P7").Value = 0.1209225617494376
$ws.Range("Q7").Value = 30.02797569783967
$ws.Range("R7").Value = 270.251781280557
$ws.Range("S7").Value = 0.1008986457783277
$ws.Range("T7").Value = 0.1008986457783277

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Slit2"
$ws.Range("C8").Value = "Gpc1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.277274333333334
$ws.Range("H8").Value = 12.831823
$ws.Range("I8").Value = 0.8344071140950421
$ws.Range("J8").Value = 0.8344071140950421
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.462291666666667
$ws.Range("N8").Value = 4.386875
$ws.Range("O8").Value = 0.02518734518879435
$ws.Range("P8").Value = 0.02518734518879435
$ws.Range("Q8").Value = 6.254622613680556
$ws.Range("R8").Value = 56.291603523125
$ws.Range("S8").Value = 0.02101650001069753
$ws.Range("T8").Value = 0.02101650001069753

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Slit2"
$ws.Range("C9").Value = "Gpc1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.277274333333334
$ws.Range("H9").Value = 12.831823
$ws.Range("I9").Value = 0.8344071140950421
$ws.Range("J9").Value = 0.8344071140950421
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 48.72612633333333
$ws.Range("N9").Value = 146.178379
$ws.Range("O9").Value = 0.8392865743864156
$ws.Range("P9").Value = 0.8392865743864156
$ws.Range("Q9").Value = 208.4150095283241
$ws.Range("R9").Value = 1875.735085754917
$ws.Range("S9").Value = 0.7003066884324829
$ws.Range("T9").Value = 0.7003066884324829

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit2"
$ws.Range("C10").Value = "Gpc1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7049536666666666
$ws.Range("H10").Value = 2.114861
$ws.Range("I10").Value = 0.1375217740863597
$ws.Range("J10").Value = 0.1375217740863597
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8478306666666667
$ws.Range("N10").Value = 2.543492
$ws.Range("O10").Value = 0.01460351867535248
$ws.Range("P10").Value = 0.01460351867535248
$ws.Range("Q10").Value = 0.5976813371791111
$ws.Range("R10").Value = 5.379132034612
$ws.Range("S10").Value = 0.002008301796137759
$ws.Range("T10").Value = 0.00200830179613776

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Slit2"
$ws.Range("C11").Value = "Gpc1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7049536666666666
$ws.Range("H11").Value = 2.114861
$ws.Range("I11").Value = 0.1375217740863597
$ws.Range("J11").Value = 0.1375217740863597
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.020353
$ws.Range("N11").Value = 21.061059
$ws.Range("O11").Value = 0.1209225617494376
$ws.Range("P11").Value = 0.1209225617494376
$ws.Range("Q11").Value = 4.949023588644333
$ws.Range("R11").Value = 44.541212297799
$ws.Range("S11").Value = 0.01662948521885003
$ws.Range("T11").Value = 0.01662948521885004

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Slit2"
$ws.Range("C12").Value = "Gpc1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7049536666666666
$ws.Range("H12").Value = 2.114861
$ws.Range("I12").Value = 0.1375217740863597
$ws.Range("J12").Value = 0.1375217740863597
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.462291666666667
$ws.Range("N12").Value = 4.386875
$ws.Range("O12").Value = 0.02518734518879435
$ws.Range("P12").Value = 0.02518734518879435
$ws.Range("Q12").Value = 1.030847872152778
$ws.Range("R12").Value = 9.277630849374999
$ws.Range("S12").Value = 0.003463808394888536
$ws.Range("T12").Value = 0.003463808394888536

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Slit2"
$ws.Range("C13").Value = "Gpc1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7049536666666666
$ws.Range("H13").Value = 2.114861
$ws.Range("I13").Value = 0.1375217740863597
$ws.Range("J13").Value = 0.1375217740863597
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 48.72612633333333
$ws.Range("N13").Value = 146.178379
$ws.Range("O13").Value = 0.8392865743864156
$ws.Range("P13").Value = 0.8392865743864156
$ws.Range("Q13").Value = 34.34966142114655
$ws.Range("R13").Value = 309.146952790319
$ws.Range("S13").Value = 0.1154201786764834
$ws.Range("T13").Value = 0.1154201786764834
